# Auto-generated edit script applying the Ravana_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 43.75
$ws.Range("I8").Value = 43.75
$ws.Range("K8").Value = 131.25
$ws.Range("M8").Value = 7.75

$ws.Range("H11").Value = 98.5
$ws.Range("I11").Value = 98.5
$ws.Range("K11").Value = 98.5
$ws.Range("M11").Value = 41.5

$ws.Range("H28").Value = 1541.5714
$ws.Range("I28").Value = 1187.6666
$ws.Range("J28").Value = 2178.6
$ws.Range("K28").Value = 1187.6666
$ws.Range("L28").Value = 2178.6
$ws.Range("M28").Value = -702.6666
$ws.Range("N28").Value = -3148.6

$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876

$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380

$ws.Range("H70").Value = 1112243.6
$ws.Range("I70").Value = 1066.3334
$ws.Range("J70").Value = 1667832.4
$ws.Range("K70").Value = 3199.0002
$ws.Range("L70").Value = 5003497.199999999
$ws.Range("M70").Value = -2929.0002
$ws.Range("N70").Value = -5004037.199999999

$ws.Range("H73").Value = 1112243.6
$ws.Range("I73").Value = 1066.3334
$ws.Range("J73").Value = 1667832.4
$ws.Range("K73").Value = 3199.0002
$ws.Range("L73").Value = 5003497.199999999
$ws.Range("M73").Value = -2263.0002
$ws.Range("N73").Value = -5005369.199999999

$ws.Range("H80").Value = 754
$ws.Range("I80").Value = 694.4
$ws.Range("K80").Value = 2083.2
$ws.Range("M80").Value = -1085.2

$ws.Range("H83").Value = 754
$ws.Range("I83").Value = 694.4
$ws.Range("K83").Value = 6249.599999999999
$ws.Range("M83").Value = -1257.599999999999

$ws.Range("H129").Value = 2583.75
$ws.Range("J129").Value = 2993
$ws.Range("L129").Value = 8979
$ws.Range("N129").Value = -18979

$ws.Range("H132").Value = 757
$ws.Range("I132").Value = 751.2143
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2253.6429
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = 276.3571000000002
$ws.Range("N132").Value = -8060

$ws.Range("H138").Value = 4373.067
$ws.Range("I138").Value = 3998.6667
$ws.Range("J138").Value = 4414.6665
$ws.Range("K138").Value = 11996.0001
$ws.Range("L138").Value = 13243.9995
$ws.Range("M138").Value = -6856.000100000001
$ws.Range("N138").Value = -23523.9995


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4858.3594
$ws.Range("I32").Value = 4682.058
$ws.Range("J32").Value = 9912.333000000001
$ws.Range("K32").Value = 4682.058
$ws.Range("L32").Value = 9912.333000000001
$ws.Range("M32").Value = -4395.058
$ws.Range("N32").Value = -10486.333

$ws.Range("H61").Value = 2076.0557
$ws.Range("I61").Value = 1527.7
$ws.Range("K61").Value = 1527.7
$ws.Range("M61").Value = -1315.7

$ws.Range("H132").Value = 2795.4092
$ws.Range("I132").Value = 1964.5714
$ws.Range("J132").Value = 4249.375
$ws.Range("K132").Value = 5893.7142
$ws.Range("L132").Value = 12748.125
$ws.Range("M132").Value = -3363.7142
$ws.Range("N132").Value = -17808.125

$ws.Range("H136").Value = 2076.0557
$ws.Range("I136").Value = 1527.7
$ws.Range("K136").Value = 4583.1
$ws.Range("M136").Value = -2033.1


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 792.625
$ws.Range("I80").Value = 1051.6
$ws.Range("K80").Value = 1051.6
$ws.Range("M80").Value = -53.59999999999991

$ws.Range("H83").Value = 792.625
$ws.Range("I83").Value = 1051.6
$ws.Range("K83").Value = 5258
$ws.Range("M83").Value = -266


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1307.8
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 1307.8
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H112").Value = 6511.706
$ws.Range("J112").Value = 6862.4375
$ws.Range("L112").Value = 20587.3125
$ws.Range("N112").Value = -22803.3125

$ws.Range("H138").Value = 7535.4
$ws.Range("I138").Value = 5892.3335
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 17677.0005
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = -12537.0005
$ws.Range("N138").Value = -40280

$ws.Range("H140").Value = 1106.7142
$ws.Range("I140").Value = 1106.7142
$ws.Range("K140").Value = 3320.1426
$ws.Range("M140").Value = 1859.8574


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J16").Value = 450
$ws.Range("L16").Value = 450
$ws.Range("N16").Value = -790

$ws.Range("H55").Value = 392.16666
$ws.Range("I55").Value = 394.6
$ws.Range("J55").Value = 380
$ws.Range("K55").Value = 394.6
$ws.Range("L55").Value = 380
$ws.Range("M55").Value = -221.6
$ws.Range("N55").Value = -726

$ws.Range("H68").Value = 1995
$ws.Range("I68").Value = 1995
$ws.Range("K68").Value = 1995
$ws.Range("M68").Value = -1246

$ws.Range("H71").Value = 1995
$ws.Range("I71").Value = 1995
$ws.Range("K71").Value = 9975
$ws.Range("M71").Value = -6231

$ws.Range("H132").Value = 5067.4
$ws.Range("I132").Value = 4910.75
$ws.Range("J132").Value = 5694
$ws.Range("K132").Value = 14732.25
$ws.Range("L132").Value = 17082
$ws.Range("M132").Value = -12202.25
$ws.Range("N132").Value = -22142

$ws.Range("H136").Value = 3560.121
$ws.Range("I136").Value = 3728.2222
$ws.Range("K136").Value = 11184.6666
$ws.Range("M136").Value = -8634.6666


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3325.8
$ws.Range("I132").Value = 2403.0908
$ws.Range("J132").Value = 4453.5557
$ws.Range("K132").Value = 7209.2724
$ws.Range("L132").Value = 13360.6671
$ws.Range("M132").Value = -4679.2724
$ws.Range("N132").Value = -18420.6671

